# Applies the "cambios en los analisis" update to the DM-test results table.
#
# What actually changed in the source data: the LSPMW model's row, within
# every fixed-d block (A column), was recomputed/re-positioned so that it now
# sits immediately after the LSPM row (sorted by the ECRPS_Sin_Diff value in
# column C, descending, same as the rest of each block). Concretely, for each
# d-group the LSPMW entry was pulled out, every row between its old slot and
# the LSPM row shifted up by one, and a (slightly different, re-computed)
# LSPMW row was inserted right after LSPM. All other rows/blocks are
# untouched. The shared-strings table also lost its separate "LSPMW" slot
# next to "EnCQR-LSTM"/"DeepAR" and gained a new one after "LSPM" as a
# consequence - Excel manages that automatically as we assign cell values
# below, so we only need to touch the actual data cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, Modelo, ECRPS_Sin_Diff, ECRPS_Con_Diff, Mejora_%, DM_Stat, p_valor
# (numbers are written in plain decimal, not E-notation - this PowerShell
# dialect's tokenizer doesn't accept scientific-notation numeric literals)
$updates = @(
    @(5, "DeepAR", 2.65747106496605, 0.5703222791555561, 78.53890916540075, 2.494871527303883, 0.02059645052341619),
    @(6, "AV-MCPS", 3.045061224774464, 0.7074335935529131, 76.76783679102182, 3.979940955020395, 0.0006333609825530662),
    @(7, "MCPS", 2.701955361016653, 0.6889183596399426, 74.50297034586364, 4.233609761091651, 0.0003410168527335777),
    @(8, "LSPM", 1.134104262201696, 0.6646327975475445, 41.39579404654932, 9.32770205465871, 0.000000004216664573064577),
    @(9, "LSPMW", 1.136768169855944, 0.671980650564856, 40.88674644628622, 9.870731880388252, 0.000000001528800641636963),
    @(13, "DeepAR", 1032.364311169534, 30.27576925938623, 97.06733670160607, 4.643670539851906, 0.0001252677480902964),
    @(15, "MCPS", 607.6788956355718, 29.35904219599299, 95.16865857826339, 4.286107696026305, 0.0002999616353533963),
    @(16, "AV-MCPS", 496.1349728514716, 29.2968893301947, 94.09497597764282, 4.505959219671162, 0.0001753018583174892),
    @(17, "LSPM", 58.16326180633573, 28.91359804756351, 50.28889861122961, 4.462958617510705, 0.0001947137330660098),
    @(18, "LSPMW", 58.13946918268717, 28.9076813004957, 50.27873197524158, 4.45838598591608, 0.0001969007947060231),
    @(20, "DeepAR", 173959.1629817165, 1018.038395667763, 99.41478311448603, 5.675215368924245, 0.00001045074495409359),
    @(21, "MCPS", 94891.83174998802, 597.7100550895532, 99.37011432484057, 5.666122478986837, 0.00001067804402676842),
    @(22, "AREPD", 224661.9279532027, 1668.793212724136, 99.25719803621033, 5.791045711614122, 0.000007951149563867332),
    @(23, "Block Bootstrapping", 239029.2143902132, 1822.124109082956, 99.23769815596334, 5.817629428354874, 0.00000746910221427477),
    @(24, "AV-MCPS", 51463.44232116333, 515.8738629160657, 98.9975916113487, 5.406969511476652, 0.0000197813312365902),
    @(25, "EnCQR-LSTM", 76725.4774647124, 804.4784241064688, 98.9514846297614, 6.144802766344924, 0.000003481112960557198),
    @(26, "LSPMW", 4687.969884107692, 54.92132823176864, 98.82846243492406, 5.762959678840502, 0.00000849495688237667),
    @(27, "LSPM", 4680.645402168596, 54.94861344143851, 98.82604622396775, 5.745753801445408, 0.000008846669779405403),
    @(29, "LSPM", 401817.5202235829, 4335.482842240131, 98.92103190527189, 6.04567001826597, 0.00000438156400450751),
    @(30, "LSPMW", 400134.7069724926, 4331.448889646987, 98.91750232754872, 6.037894737575589, 0.000004461552906676047),
    @(31, "MCPS", 7650016.09706249, 85313.40533540784, 98.88479443372457, 6.046966344571899, 0.000004368371011587158),
    @(32, "DeepAR", 14005904.0100794, 160237.4993842074, 98.8559289049183, 5.720821364936925, 0.00000938282756179909),
    @(34, "AV-MCPS", 3949555.41930737, 50728.35286098456, 98.71559333961997, 6.255841783842101, 0.00000269398390417841),
    @(38, "LSPM", 31959278.30984279, 482293.4133794098, 98.49091268988114, 5.858055782967347, 0.000006792433822955246),
    @(39, "LSPMW", 31383316.41009646, 479389.9575332087, 98.47246877522802, 5.975061803818217, 0.00000516519300819418),
    @(40, "AV-MCPS", 302029135.8022056, 4652938.886985654, 98.45944038656165, 5.99361615974306, 0.000004946358789670313),
    @(42, "MCPS", 574225765.8906945, 9165618.36345561, 98.40383018180337, 6.011397089147072, 0.000004745522187876716),
    @(43, "DeepAR", 992793839.8827198, 16875825.74089928, 98.30016816553849, 6.047728025981449, 0.000004360638136802208),
    @(48, "LSPM", 1287459788.280791, 25243545.10549158, 98.03927506433419, 7.965213233173771, 0.00000006345664749218827),
    @(49, "LSPMW", 1265644238.911851, 24942430.48433685, 98.02927001779099, 8.033005460928365, 0.00000005512877110369629),
    @(50, "AV-MCPS", 12023767741.14577, 241300287.4162372, 97.99313915062997, 8.04942311779922, 0.00000005328705965368385),
    @(51, "MCPS", 22120957776.42364, 455451865.039818, 97.94108433439878, 8.007952831548478, 0.00000005806658842644197),
    @(53, "DeepAR", 36559212319.26645, 796011138.9636092, 97.82267973387349, 7.706918887139473, 0.0000001090728460884094),
    @(57, "LSPM", 64380602828.17851, 1541935848.40369, 97.60496829686596, 5.706611632530271, 0.000009703087466617077),
    @(58, "LSPMW", 62460288396.99023, 1523777560.795277, 97.56040581959799, 5.786140801158132, 0.000008043497842891156),
    @(59, "AV-MCPS", 592496493073.3569, 14711932116.51295, 97.51695878566298, 5.713297110069927, 0.000009551046070699343),
    @(61, "MCPS", 1052859501170.646, 27058168820.91214, 97.43003042753313, 5.746966964764621, 0.000008821391913160781),
    @(62, "DeepAR", 1620746478918.832, 43858732466.16466, 97.29391776958097, 5.799533888521103, 0.00000779388080918153),
    @(65, "LSPM", 1476931647182352, 54869063329978.05, 96.28492872810629, 7.165013066362093, 0.0000003498223013487234),
    @(66, "LSPMW", 1436397544180646, 53371990299098.23, 96.28431623854225, 7.226498109201343, 0.0000003058914341824703),
    @(67, "AV-MCPS", 12732581342679240, 486637300837981.1, 96.17801537849368, 7.272712060199441, 0.0000002766339517457084),
    @(69, "MCPS", 20794986822359090, 817390885859121.2, 96.0692887529015, 7.25010344191222, 0.0000002905696052657447),
    @(70, "DeepAR", 28702183129570310, 1153203812976251, 95.98217387238336, 7.203584202660262, 0.0000003215595709704644)
)

foreach ($u in $updates) {
    $r = $u[0]
    $ws.Cells.Item($r, 2).Value = $u[1]
    $ws.Cells.Item($r, 3).Value = $u[2]
    $ws.Cells.Item($r, 4).Value = $u[3]
    $ws.Cells.Item($r, 5).Value = $u[4]
    $ws.Cells.Item($r, 6).Value = $u[5]
    $ws.Cells.Item($r, 7).Value = $u[6]
}
